$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 19374.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 19374.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 19374.5
$ws.Range("N108").Value = -27054.5

# Hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2000
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -4253
$ws.Range("N121").ClearContents()

# Hunk 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4790.645
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 4937
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 14811
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -24891

# Hunk 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1607.174
$ws.Range("I141").Value = 1203.421
$ws.Range("J141").Value = 3525
$ws.Range("K141").Value = 3610.263
$ws.Range("L141").Value = 10575
$ws.Range("M141").Value = 1569.737
$ws.Range("N141").Value = -20935

# Hunk 4: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 85519.086
$ws.Range("I2").Value = 2489.3333
$ws.Range("J2").Value = 168548.83
$ws.Range("K2").Value = 2489.3333
$ws.Range("L2").Value = 168548.83
$ws.Range("M2").Value = -2376.3333
$ws.Range("N2").Value = -168774.83

# Hunk 5: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25473.549
$ws.Range("I32").Value = 4627.857
$ws.Range("J32").Value = 156801.4
$ws.Range("K32").Value = 4627.857
$ws.Range("L32").Value = 156801.4
$ws.Range("M32").Value = -4340.857

# Hunk 6: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 58001.723
$ws.Range("I45").Value = 145204.28
$ws.Range("J45").Value = 2509.182
$ws.Range("K45").Value = 145204.28
$ws.Range("L45").Value = 2509.182
$ws.Range("M45").Value = -144827.28

# Hunk 7: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 9499.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 9499.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 9499.5
$ws.Range("N98").Value = -15489.5

# Hunk 8: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 85519.086
$ws.Range("I116").Value = 2489.3333
$ws.Range("J116").Value = 168548.83
$ws.Range("K116").Value = 2489.3333
$ws.Range("L116").Value = 168548.83
$ws.Range("M116").Value = -195.3332999999998
$ws.Range("N116").Value = -173136.83

# Hunk 9: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Hunk 10: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2335.2666
$ws.Range("I132").Value = 1926.3784
$ws.Range("J132").Value = 4226.375
$ws.Range("K132").Value = 5779.135200000001
$ws.Range("L132").Value = 12679.125
$ws.Range("M132").Value = -3249.135200000001
$ws.Range("N132").Value = -17739.125

# Hunk 11: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 67200
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 67200
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 67200
$ws.Range("N140").Value = -77560

# Hunk 12: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 85519.086
$ws.Range("I3").Value = 2489.3333
$ws.Range("J3").Value = 168548.83
$ws.Range("K3").Value = 2489.3333
$ws.Range("L3").Value = 168548.83
$ws.Range("M3").Value = -2375.3333
$ws.Range("N3").Value = -168776.83

# Hunk 13: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47246.047
$ws.Range("I20").Value = 60707.234
$ws.Range("J20").Value = 1478
$ws.Range("K20").Value = 60707.234
$ws.Range("L20").Value = 1478
$ws.Range("M20").Value = -60460.234
$ws.Range("N20").Value = -1972

# Hunk 14: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()

# Hunk 15: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 19950
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 19950
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 19950
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -20298

# Hunk 16: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26255.834
$ws.Range("I31").Value = 42947.207
$ws.Range("J31").Value = 4000.6667
$ws.Range("K31").Value = 42947.207
$ws.Range("L31").Value = 4000.6667
$ws.Range("M31").Value = -42652.207
$ws.Range("N31").Value = -4590.6667

# Hunk 17: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 26255.834
$ws.Range("I34").Value = 42947.207
$ws.Range("J34").Value = 4000.6667
$ws.Range("K34").Value = 42947.207
$ws.Range("L34").Value = 4000.6667
$ws.Range("M34").Value = -42745.207
$ws.Range("N34").Value = -4404.6667

# Hunk 18: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6946719
$ws.Range("I62").Value = 27778780
$ws.Range("J62").Value = 2698.6667
$ws.Range("K62").Value = 27778780
$ws.Range("L62").Value = 2698.6667
$ws.Range("M62").Value = -27778156
$ws.Range("N62").Value = -3946.6667

# Hunk 19: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 6946719
$ws.Range("I65").Value = 27778780
$ws.Range("J65").Value = 2698.6667
$ws.Range("K65").Value = 138893900
$ws.Range("L65").Value = 13493.3335
$ws.Range("M65").Value = -138890780
$ws.Range("N65").Value = -19733.3335

# Hunk 20: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14257.667
$ws.Range("I99").Value = 3548.5
$ws.Range("J99").Value = 35676
$ws.Range("K99").Value = 3548.5
$ws.Range("L99").Value = 35676
$ws.Range("M99").Value = -2050.5

# Hunk 21: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 14257.667
$ws.Range("I126").Value = 3548.5
$ws.Range("J126").Value = 35676
$ws.Range("K126").Value = 10645.5
$ws.Range("L126").Value = 107028
$ws.Range("M126").Value = -8175.5

# Hunk 22: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6757.758
$ws.Range("I5").Value = 592.04346
$ws.Range("J5").Value = 20938.9
$ws.Range("K5").Value = 1776.13038
$ws.Range("L5").Value = 62816.7
$ws.Range("M5").Value = -1664.13038

# Hunk 23: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 464.33334
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 464.33334
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1393.00002
$ws.Range("N92").Value = -3889.00002

# Hunk 24: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 619.3043
$ws.Range("I113").Value = 580
$ws.Range("J113").Value = 640.26666
$ws.Range("K113").Value = 1740
$ws.Range("L113").Value = 1920.79998
$ws.Range("M113").Value = 430
$ws.Range("N113").Value = -6260.79998

# Hunk 25: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4634.0835
$ws.Range("I122").Value = 361
$ws.Range("J122").Value = 25999.5
$ws.Range("K122").Value = 3249
$ws.Range("L122").Value = 233995.5
$ws.Range("M122").Value = -799
$ws.Range("N122").Value = -238895.5

# Hunk 26: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 750.01
$ws.Range("I131").Value = 468.2143
$ws.Range("J131").Value = 795.8837
$ws.Range("K131").Value = 1404.6429
$ws.Range("L131").Value = 2387.6511
$ws.Range("M131").Value = 3635.3571
$ws.Range("N131").Value = -12467.6511

# Hunk 27: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6757.758
$ws.Range("I135").Value = 592.04346
$ws.Range("J135").Value = 20938.9
$ws.Range("K135").Value = 5328.39114
$ws.Range("L135").Value = 188450.1
$ws.Range("M135").Value = -2793.39114

# Hunk 28: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2030336
$ws.Range("I126").Value = 1714.2307
$ws.Range("J126").Value = 3678591.2
$ws.Range("K126").Value = 5142.6921
$ws.Range("L126").Value = 11035773.6
$ws.Range("M126").Value = -2672.6921
$ws.Range("N126").Value = -11040713.6

# Hunk 29: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 33075.332
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 33075.332
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 99225.99600000001
$ws.Range("N134").Value = -104295.996

# Hunk 30: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3020
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 4550
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 4550
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -4774

# Hunk 31: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 39263.58
$ws.Range("I16").Value = 62853.125
$ws.Range("J16").Value = 1520.3
$ws.Range("K16").Value = 62853.125
$ws.Range("L16").Value = 1520.3
$ws.Range("M16").Value = -62683.125
$ws.Range("N16").Value = -1860.3

# Hunk 32: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 500999.5
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 1999
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 1999
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -2271

# Hunk 33: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3020
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4550
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 13650
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -18590

# Hunk 34: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 40950
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 40950
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 40950
$ws.Range("N128").Value = -50910

# Hunk 35: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5559595.5
$ws.Range("I62").Value = 27778378
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 27778378
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -27777754
$ws.Range("N62").Value = -6148

# Hunk 36: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5559595.5
$ws.Range("I65").Value = 27778378
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 138891890
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -138888770
$ws.Range("N65").Value = -30740

# Hunk 37: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 34495
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34495
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34495
$ws.Range("N124").Value = -44315

# Hunk 38: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2536.8
$ws.Range("I126").Value = 2536.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7610.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5140.400000000001
$ws.Range("N126").ClearContents()
